# Add an "alpha" column (1s) into the "17.02.2015" sheet, inserted before
# the existing "trained on gold / high order opt" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("17.02.2015")
$ws.Activate()

# Insert a new column at J, shifting existing J:K -> K:L. The new column
# inherits the width of its left neighbour (I), matching Excel's default
# insert behaviour.
$leftWidth = $ws.Columns("I").ColumnWidth
$ws.Columns("J:J").Insert()
$ws.Columns("J:J").ColumnWidth = $leftWidth

# Header for the new column.
$ws.Range("J1").Value = "alpha"

# Fill rows 2-21 with the constant 1.
$ws.Range("J2:J21").Value = 1

# Re-establish the AutoFilter over the (now 12-column) table range A1:L13.
# Calling AutoFilter on a range whose rows are contiguous with more data
# below it snaps to the full contiguous block, so the rows below the table
# (14:21) are stashed and blanked first, restored once the filter is set.
$saved = $ws.Range("A14:L21").Value()
$ws.Range("A14:L21").ClearContents()
$ws.AutoFilterMode = $false
$ws.Range("A1:L13").AutoFilter() | Out-Null
$ws.Range("A14:L21").Value = $saved

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$fdName = $wb.Names.Item("17.02.2015!_FilterDatabase")
$fdName.RefersTo = "='17.02.2015'!`$A`$1:`$L`$13"

# Match the recorded selection after the edit.
$ws.Range("J7").Select()
